$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "2 mins"
$ws.Range("C3").Value = "1 min 47 seconds"
$ws.Range("D3").Value = "have just used the static method from ItemDataLoader"
$ws.Range("E3").Value = "I may use the public static methods from ItemDataLoader, which could make tasks easier"

$ws.Rows.Item(3).RowHeight = 45

$ws.Range("B4:E8").Select()
